# Updates cryptos list values (Price / Volume(1h) columns, plus a row swap)
# as captured in the "Updated cryptos list" GitHub Actions commit.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Helper: write a value as plain text, preserving the default (unstyled) cell
# format even when the text looks like a number (e.g. "372.21" or "1.00").
# Excel normally auto-converts such strings to numeric values; forcing the
# NumberFormat to Text ("@") before the assignment keeps it a string, and
# resetting the Style back to "Normal" afterwards avoids leaving a stray
# number-format style applied to the cell.
function Set-TextValue($cell, $value) {
    $c = $ws.Range($cell)
    $c.NumberFormat = "@"
    $c.Value = $value
    $c.Style = "Normal"
}

$ws.Range('D2').Value = '51.310.37'
$ws.Range('E2').Value = '  -1.78%  '
$ws.Range('D3').Value = '2.922.10'
$ws.Range('E3').Value = '  -0.64%  '
$ws.Range('E4').Value = '  -0.29%  '
Set-TextValue 'D5' '372.21'
$ws.Range('E5').Value = '  +4.33%  '
Set-TextValue 'D6' '103.84'
$ws.Range('E6').Value = '  -5.22%  '
Set-TextValue 'D7' '0.541'
$ws.Range('E7').Value = '  -5.13%  '
Set-TextValue 'D8' '1.00'
Set-TextValue 'D9' '0.592'
$ws.Range('E9').Value = '  -5.82%  '
Set-TextValue 'D10' '37.19'
$ws.Range('E10').Value = '  -4.80%  '
$ws.Range('E11').Value = '  +0.76%  '
Set-TextValue 'D12' '0.0839'
$ws.Range('E12').Value = '  -4.25%  '
Set-TextValue 'D13' '18.46'
$ws.Range('E13').Value = '  -5.77%  '
$ws.Range('D14').Value = '3.378.40'
$ws.Range('E14').Value = '  -0.92%  '
Set-TextValue 'D15' '7.40'
$ws.Range('E15').Value = '  -5.42%  '
$ws.Range('D16').Value = '2.906.20'
$ws.Range('E16').Value = '  -1.20%  '
Set-TextValue 'D17' '0.953'
$ws.Range('E17').Value = '  -2.88%  '
$ws.Range('D18').Value = '51.221.18'
$ws.Range('E18').Value = '  -1.97%  '
Set-TextValue 'D19' '3.34'
$ws.Range('E19').Value = '  -5.84%  '
Set-TextValue 'D20' '7.31'
$ws.Range('E20').Value = '  -3.90%  '
Set-TextValue 'D21' '13.09'
$ws.Range('E21').Value = '  -5.94%  '
$ws.Range('D22').Value = '0.0₃0948'
$ws.Range('E22').Value = '  -3.47%  '
Set-TextValue 'D23' '68.38'
$ws.Range('E23').Value = '  -3.04%  '
Set-TextValue 'D24' '261.93'
$ws.Range('E24').Value = '  -3.42%  '
Set-TextValue 'D25' '2.72'
$ws.Range('E25').Value = '  -2.74%  '
Set-TextValue 'D26' '4.35'
$ws.Range('E26').Value = '  +3.98%  '
$ws.Range('E27').Value = '  -3.06%  '
$ws.Range('E28').Value = '  +0.10%  '
$ws.Range('E29').Value = '  -6.06%  '
Set-TextValue 'D30' '25.99'
$ws.Range('E30').Value = '  -3.73%  '
Set-TextValue 'D31' '0.105'
$ws.Range('E31').Value = '  -3.84%  '
$ws.Range('B32').Value = 'Cosmos'
$ws.Range('C32').Value = 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom'
Set-TextValue 'D32' '9.95'
$ws.Range('E32').Value = '  -5.14%  '
$ws.Range('B33').Value = 'RenderToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
Set-TextValue 'D33' '6.14'
$ws.Range('E33').Value = '  -1.34%  '
$ws.Range('E34').Value = '  -6.35%  '
Set-TextValue 'D35' '35.24'
$ws.Range('E35').Value = '  -6.38%  '
Set-TextValue 'D36' '50.73'
$ws.Range('E36').Value = '  -2.60%  '
$ws.Range('E37').Value = '  +0.03%  '
Set-TextValue 'D38' '0.0426'
$ws.Range('E38').Value = '  -4.20%  '
Set-TextValue 'D39' '2.77'
$ws.Range('E39').Value = '  +0.63%  '
$ws.Range('E40').Value = '  -1.86%  '
Set-TextValue 'D41' '17.10'
$ws.Range('E41').Value = '  -6.53%  '
Set-TextValue 'D42' '1.87'
$ws.Range('E42').Value = '  -6.45%  '
$ws.Range('E43').Value = '  -5.62%  '
Set-TextValue 'D44' '22.22'
$ws.Range('E44').Value = '  -3.09%  '
Set-TextValue 'D45' '117.43'
$ws.Range('E45').Value = '  -1.90%  '
Set-TextValue 'D46' '2.09'
$ws.Range('E46').Value = '  -3.65%  '
$ws.Range('D47').Value = '2.060.78'
$ws.Range('E47').Value = '  -3.41%  '
Set-TextValue 'D48' '2.31'
$ws.Range('E48').Value = '  -6.35%  '
Set-TextValue 'D49' '3.21'
$ws.Range('E49').Value = '  -7.48%  '
$ws.Range('D50').Value = '3.219.52'
$ws.Range('E50').Value = '  -0.46%  '
Set-TextValue 'D51' '0.235'
$ws.Range('E51').Value = '  -5.73%  '
